$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-12 Thursday" "2026-02-13 Friday"

Replace-Text "304×3=912" "572×4=2288"
Replace-Text "434×5=2170" "537×5=2685"
Replace-Text "343×4=1372" "416×6=2496"
Replace-Text "358×9=3222" "726×5=3630"
Replace-Text "849×8=6792" "262×7=1834"

Replace-Text "591×9=5319" "867×7=6069"
Replace-Text "899×6=5394" "378×3=1134"
Replace-Text "207×8=1656" "447×2=894"
Replace-Text "333×8=2664" "281×4=1124"
Replace-Text "744×5=3720" "574×4=2296"

Replace-Text "862×8=6896" "858×7=6006"
Replace-Text "558×2=1116" "808×8=6464"
Replace-Text "492×9=4428" "624×4=2496"
Replace-Text "243×9=2187" "637×7=4459"
Replace-Text "940×3=2820" "945×4=3780"

Replace-Text "440×3=1320" "766×6=4596"
Replace-Text "516×7=3612" "367×7=2569"
Replace-Text "236×5=1180" "706×8=5648"
Replace-Text "426×8=3408" "776×3=2328"
Replace-Text "841×8=6728" "550×5=2750"

Replace-Text "745×3=2235" "358×8=2864"
Replace-Text "330×6=1980" "702×6=4212"
Replace-Text "791×6=4746" "182×4=728"
Replace-Text "650×4=2600" "754×5=3770"
Replace-Text "502×9=4518" "644×9=5796"
